# HeroPower.xlsx / HeroSkill sheet: add a new "skill11" equip/power row to
# the "表1" table (shield power - grants target unit some magic armor).
#
# add some new equips. remove some useless module

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the table ("表1") by one row - ListRows.Add() extends the table
# range/autofilter and keeps the table's banded-row formatting in sync,
# same as typing a new row right under the table in the UI.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# Copy the formatting of the previous last row down onto the freshly
# inserted row so the new cells pick up style "s=1" like every other
# data row instead of staying unformatted.
$ws.Range("A13:F13").Copy()
$ws.Range("A14:F14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New row: Id, Name, Des, Type, CardId, Icon
$ws.Range("A14").Value = 31100001
$ws.Range("B14").Value = "护盾"
$ws.Range("C14").Value = "使指定单位获得一些魔甲"
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 53200001
$ws.Range("F14").Value = "skill11"

# Match the author's final selection/cursor position.
[void]$ws.Range("F14").Select()
